# Hortaliza, Vega Central Mapocho de Santiago - Ajo
# Insert 3 new weekly price rows (Ajo, Rosado, "nueva(o)" grades) before the
# existing row 74, shifting all subsequent rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows starting at row 74 (pushes old rows 74:151 -> 77:154)
$ws.Rows("74:76").Insert()

# Values shared by every data row in this sheet (constant columns)
$marketId = 9
$market   = "Vega Central Mapocho de Santiago"
$region   = "Metropolitana"
$codreg   = 13
$catId    = 100112003
$category = "Ajo"
$clasif   = "Hortaliza"

# New row 74: Ajo / Rosado / 1a nueva(o)
$r = 74
$ws.Cells.Item($r, 1).Value  = $marketId
$ws.Cells.Item($r, 2).Value  = $market
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 44539
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $catId
$ws.Cells.Item($r, 7).Value  = $category
$ws.Cells.Item($r, 8).Value  = "Rosado"
$ws.Cells.Item($r, 9).Value  = "1a nueva(o)"
$ws.Cells.Item($r, 10).Value = 5000
$ws.Cells.Item($r, 11).Value = 1800
$ws.Cells.Item($r, 12).Value = 2000
$ws.Cells.Item($r, 13).Value = 1900
$ws.Cells.Item($r, 14).Value = "$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item($r, 15).Value = "Provincia de Talagante"
$ws.Cells.Item($r, 16).Value = 95
$ws.Cells.Item($r, 17).Value = 20
$ws.Cells.Item($r, 18).Value = $clasif

# New row 75: Ajo / Rosado / 2a nueva(o)
$r = 75
$ws.Cells.Item($r, 1).Value  = $marketId
$ws.Cells.Item($r, 2).Value  = $market
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 44539
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $catId
$ws.Cells.Item($r, 7).Value  = $category
$ws.Cells.Item($r, 8).Value  = "Rosado"
$ws.Cells.Item($r, 9).Value  = "2a nueva(o)"
$ws.Cells.Item($r, 10).Value = 3200
$ws.Cells.Item($r, 11).Value = 1200
$ws.Cells.Item($r, 12).Value = 1500
$ws.Cells.Item($r, 13).Value = 1350
$ws.Cells.Item($r, 14).Value = "$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item($r, 15).Value = "Provincia de Talagante"
$ws.Cells.Item($r, 16).Value = 68
$ws.Cells.Item($r, 17).Value = 20
$ws.Cells.Item($r, 18).Value = $clasif

# New row 76: Ajo / Rosado / Extra nueva (o)
$r = 76
$ws.Cells.Item($r, 1).Value  = $marketId
$ws.Cells.Item($r, 2).Value  = $market
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 44539
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $catId
$ws.Cells.Item($r, 7).Value  = $category
$ws.Cells.Item($r, 8).Value  = "Rosado"
$ws.Cells.Item($r, 9).Value  = "Extra nueva (o)"
$ws.Cells.Item($r, 10).Value = 1400
$ws.Cells.Item($r, 11).Value = 2400
$ws.Cells.Item($r, 12).Value = 2500
$ws.Cells.Item($r, 13).Value = 2450
$ws.Cells.Item($r, 14).Value = "$/paquete 20 unidades (volumen en unidades)"
$ws.Cells.Item($r, 15).Value = "Provincia de Talagante"
$ws.Cells.Item($r, 16).Value = 122
$ws.Cells.Item($r, 17).Value = 20
$ws.Cells.Item($r, 18).Value = $clasif
